$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 51 values, mirroring the previous "no match" rows.
# Force A51 to be stored as plain text (not auto-converted to a date serial)
# the same way the other date cells in column A already are, then drop the
# number-format override again so no extra style gets attached to the cell.
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = "2025-04-16"
$ws.Range("A51").Style = "Normal"

$ws.Range("B51").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C51").Value = "NA"
$ws.Range("D51").Value = 1

# Clear the "NA" text now left behind in C50 (it moved down to C51)
$ws.Range("C50").Value = ""
